$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates a single cell to its new text value, preserving the
# original "text" storage (e.g. "169.00", "0.0000275") instead of letting
# Excel coerce it to a number and drop formatting/precision.
$updates = @(
    @{ Cell = 'D2'; Value = '67.404.22' }
    @{ Cell = 'E2'; Value = '  -1.23%  ' }
    @{ Cell = 'D3'; Value = '3.752.63' }
    @{ Cell = 'E3'; Value = '  -2.03%  ' }
    @{ Cell = 'E4'; Value = '  +0.10%  ' }
    @{ Cell = 'D5'; Value = '595.12' }
    @{ Cell = 'E5'; Value = '  -1.01%  ' }
    @{ Cell = 'D6'; Value = '169.00' }
    @{ Cell = 'E6'; Value = '  -0.43%  ' }
    @{ Cell = 'D7'; Value = '3.749.51' }
    @{ Cell = 'E7'; Value = '  -2.12%  ' }
    @{ Cell = 'E8'; Value = '  -0.03%  ' }
    @{ Cell = 'D9'; Value = '0.523' }
    @{ Cell = 'E9'; Value = '  -0.74%  ' }
    @{ Cell = 'E10'; Value = '  +0.18%  ' }
    @{ Cell = 'E12'; Value = '  -1.19%  ' }
    @{ Cell = 'D13'; Value = '0.0000275' }
    @{ Cell = 'E13'; Value = '  +3.56%  ' }
    @{ Cell = 'D14'; Value = '36.40' }
    @{ Cell = 'E14'; Value = '  -1.98%  ' }
    @{ Cell = 'D15'; Value = '4.383.11' }
    @{ Cell = 'E15'; Value = '  -2.01%  ' }
    @{ Cell = 'D16'; Value = '3.732.40' }
    @{ Cell = 'E16'; Value = '  -2.48%  ' }
    @{ Cell = 'D17'; Value = '18.54' }
    @{ Cell = 'E17'; Value = '  +0.40%  ' }
    @{ Cell = 'D18'; Value = '67.350.96' }
    @{ Cell = 'E18'; Value = '  -1.37%  ' }
    @{ Cell = 'D19'; Value = '7.17' }
    @{ Cell = 'E19'; Value = '  -3.10%  ' }
    @{ Cell = 'D21'; Value = '10.50' }
    @{ Cell = 'E21'; Value = '  -5.25%  ' }
    @{ Cell = 'D22'; Value = '466.04' }
    @{ Cell = 'E22'; Value = '  -0.50%  ' }
    @{ Cell = 'E23'; Value = '  -2.31%  ' }
    @{ Cell = 'B24'; Value = 'Litecoin' }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' }
    @{ Cell = 'D24'; Value = '83.58' }
    @{ Cell = 'E24'; Value = '  +0.52%  ' }
    @{ Cell = 'B25'; Value = 'PEPE' }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe' }
    @{ Cell = 'D25'; Value = '0.0000146' }
    @{ Cell = 'E25'; Value = '  -8.80%  ' }
    @{ Cell = 'E26'; Value = '  -1.65%  ' }
    @{ Cell = 'D27'; Value = '12.14' }
    @{ Cell = 'E27'; Value = '  -0.53%  ' }
    @{ Cell = 'D28'; Value = '10.27' }
    @{ Cell = 'E28'; Value = '  +2.29%  ' }
    @{ Cell = 'E29'; Value = '  +0.06%  ' }
    @{ Cell = 'E30'; Value = '  -2.36%  ' }
    @{ Cell = 'D31'; Value = '3.902.19' }
    @{ Cell = 'E31'; Value = '  -1.95%  ' }
    @{ Cell = 'D32'; Value = '7.64' }
    @{ Cell = 'E32'; Value = '  -0.84%  ' }
    @{ Cell = 'D33'; Value = '30.43' }
    @{ Cell = 'E33'; Value = '  -3.71%  ' }
    @{ Cell = 'D34'; Value = '2.22' }
    @{ Cell = 'E34'; Value = '  -3.91%  ' }
    @{ Cell = 'D35'; Value = '9.12' }
    @{ Cell = 'E35'; Value = '  -2.90%  ' }
    @{ Cell = 'D36'; Value = '3.711.42' }
    @{ Cell = 'E36'; Value = '  -2.20%  ' }
    @{ Cell = 'D37'; Value = '3.80' }
    @{ Cell = 'E37'; Value = '  +2.52%  ' }
    @{ Cell = 'E38'; Value = '  -1.48%  ' }
    @{ Cell = 'E39'; Value = '  -2.02%  ' }
    @{ Cell = 'D40'; Value = '0.995' }
    @{ Cell = 'E40'; Value = '  -2.16%  ' }
    @{ Cell = 'D41'; Value = '5.81' }
    @{ Cell = 'E41'; Value = '  -2.31%  ' }
    @{ Cell = 'D42'; Value = '1.00' }
    @{ Cell = 'E42'; Value = '  +0.04%  ' }
    @{ Cell = 'D43'; Value = '0.311' }
    @{ Cell = 'E43'; Value = '  -1.01%  ' }
    @{ Cell = 'D45'; Value = '8.69' }
    @{ Cell = 'E45'; Value = '  -0.76%  ' }
    @{ Cell = 'D46'; Value = '1.94' }
    @{ Cell = 'E46'; Value = '  -2.10%  ' }
    @{ Cell = 'D47'; Value = '45.86' }
    @{ Cell = 'E47'; Value = '  -2.70%  ' }
    @{ Cell = 'D48'; Value = '396.79' }
    @{ Cell = 'E48'; Value = '  -4.85%  ' }
    @{ Cell = 'E49'; Value = '  -7.81%  ' }
    @{ Cell = 'E50'; Value = '  -1.95%  ' }
    @{ Cell = 'D51'; Value = '138.67' }
    @{ Cell = 'E51'; Value = '  -2.37%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
